$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New parametric-simulation results for the surviving rows (2-5): columns B..K.
$data = @(
    @(0.1852904725172569, 8.115593631060179, 3.346357535337802, 13.93863149905029, 31.35664327227182, 0, 84919097014.25925, 94296351698.93303, 0, $true),
    @(0.6430740511735359, 21.0778076648316, 2.332619033092524, 10.18090776121786, 31.90433105341481, 56166773677.43803, 278342036.7269922, 62376119005.14106, 0, $true),
    @(0.7094288550566084, 18.4099040190732, 3.036001359008221, 12.20420376626922, 32.86446845411987, 177641522987.82, 440400452.1712334, 199608669951.4375, 0, $true),
    @(0.8680082502510337, 22.37146375150644, 2.833801455947056, 13.00726277875406, 30.29909187131385, 132483495389.3125, 2022039.508653594, 160590908914.5361, 0, $true)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

# Rows 6-11 are no longer part of the parametric sweep; remove them entirely
# so the used range (and dimension) shrinks back down to A1:K5.
$ws.Range("A6:K11").Delete() | Out-Null
